# Add a new "Yearly demand" worksheet at the end of the workbook and
# populate it with the hourly demand-offset table (3 profile rows x 24 hours).

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet ("Connected Households").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Yearly demand"

# Header row: 0..23 across B1:Y1
$hours = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23)
for ($i = 0; $i -lt $hours.Count; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $hours[$i]
}

# Data rows: col A holds the profile index (0,1,2); B..Y hold the values.
$rows = @(
    @(0, -32.5, -19.5, -13, -13, -13, 142.5, 291.5, 327, 388.5, 502, 596, 670.5, 745, 651, 576.5, 502, 320.5, 139, 32, -117, -97.5, -78, -52, -39),
    @(1, -32.5, -19.5, -13, 0, 0, -19.5, 0, 324, 486, 648, 729, 751.5, 583, 567, 333.5, 340, 243, 57.99999999999999, -130, 0, 0, -78, 0, -39),
    @(2, -32.5, -19.5, 0, 0, 0, -19.5, 0, 0, 81, 324, 567, 589.5, 648, 567, 324, 162, 81, 0, -130, 0, 0, 0, 0, -39)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowData = $rows[$r]
    $excelRow = 2 + $r
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($excelRow, 1 + $c).Value = $rowData[$c]
    }
}

# Match the header/index-column styling used on the other sheets (bold,
# centered, thin border all round) — this is cell style "s=1" in the model.
$headerRange = $ws.Range("B1:Y1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$indexRange = $ws.Range("A2:A4")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1
$indexRange.Borders.Weight = 2

$ws.Range("A1").Select()
